# "just some delete new item fix"
# - Row 3 (item WX001 / "WX001-Juguete") is removed entirely.
# - Row 2 (previously item WR005 / "WR005-Plastilina") is updated to a new
#   item: "ST012-Asd" / "ST012", and its sale price (column G) changes to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# First delete the whole third row (the item being removed), shifting rows
# up so nothing is left dangling below the remaining item.
$ws.Rows.Item(3).Delete()

# Update the remaining item's data (row 2) to the new/fixed item values.
$ws.Range("A2").Value = "ST012-Asd"
$ws.Range("B2").Value = "ST012"
$ws.Range("G2").Value = 1
$ws.Range("P2").Value = "ST012-Asd"
$ws.Range("T2").Value = "ST012"
